$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metadata")

# Insert two new rows before the old row 4 (dataset.commit.id), pushing
# everything else down by two rows.
$ws.Rows("4:5").Insert()

# Row 4: dataset.preview.table
$ws.Range("A4").Value = "dataset.preview.table"
$ws.Range("B4").Value = "source(ds:'{{dataset.id}}');`nquery([`n  { dim:'time', role:'row', items:[] },`n  { dim:'indicator', role:'col', items:[] } `n]);`nformat(p:3);`norder(dir:'row', index:-1, asc:'az');`nlimit(start:0, length:5);"

# Row 5: dataset.preview.line
$ws.Range("A5").Value = "dataset.preview.line"
$ws.Range("B5").Value = "source(ds:'{{dataset.id}}');`nquery([`n  { dim:'time', role:'row', items:[] },`n  { dim:'indicator', role:'col', items:[] } `n]);`nformat(p:3);`norder(dir:'row', index:-1, asc:'az');`nline(x:-1);"

# New rows get a taller height to show the multi-line preview formulas, and
# the text should wrap within the cell.
$ws.Rows("4:5").RowHeight = 120
$ws.Range("A4:B5").WrapText = $true
$ws.Range("A4:B5").VerticalAlignment = -4108

# Move the active selection to B10 (dataset.label's value cell after the shift).
$ws.Range("B10").Select() | Out-Null
